$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the C20/D20 descriptions: "polishing (db layer)" moves to C20,
# and the time range for that day gets its end time filled in (D20).
$ws.Range("C20").Value = "polishing (db layer)"
$ws.Range("D20").Value = "14:40-22:30; 23:30-02:15"

# Add a new row (21) for the unit tests work entry.
$ws.Range("C21").Value = "unit tests"
$ws.Range("D21").Value = "10:30-12:00; 13:30-zeit"

# Copy the time-format style from D20 onto the new D21 cell.
$ws.Range("D20").Copy()
$ws.Range("D21").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Update the selection to match the new active cell shown in the saved file.
$ws.Range("D21").Select()
